$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.423.21"
$ws.Range("E2").Value = "  -2.38%  "
$ws.Range("D3").Value = "2.970.27"
$ws.Range("E3").Value = "  -3.55%  "
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").Value = "'529.33"
$ws.Range("E5").Value = "  -1.66%  "
$ws.Range("D6").Value = "'129.65"
$ws.Range("E6").Value = "  -2.92%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "2.967.23"
$ws.Range("E8").Value = "  -3.29%  "
$ws.Range("D9").Value = "'0.484"
$ws.Range("E9").Value = "  -0.74%  "
$ws.Range("D10").Value = "'0.148"
$ws.Range("E10").Value = "  -2.95%  "
$ws.Range("D11").Value = "'6.05"
$ws.Range("E11").Value = "  -1.36%  "
$ws.Range("D12").Value = "'0.437"
$ws.Range("E12").Value = "  -4.46%  "
$ws.Range("D13").Value = "'0.0000215"
$ws.Range("E13").Value = "  -3.37%  "
$ws.Range("D14").Value = "'32.96"
$ws.Range("E14").Value = "  -4.20%  "
$ws.Range("D15").Value = "3.477.97"
$ws.Range("E15").Value = "  -1.72%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "61.707.57"
$ws.Range("E16").Value = "  -1.96%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "'0.110"
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("D18").Value = "3.010.95"
$ws.Range("E18").Value = "  -2.29%  "
$ws.Range("D19").Value = "'6.42"
$ws.Range("E19").Value = "  -2.44%  "
$ws.Range("D20").Value = "'454.03"
$ws.Range("E20").Value = "  -5.21%  "
$ws.Range("D21").Value = "'12.99"
$ws.Range("E21").Value = "  -2.23%  "
$ws.Range("D22").Value = "'0.666"
$ws.Range("E22").Value = "  -4.64%  "
$ws.Range("D23").Value = "'6.76"
$ws.Range("E23").Value = "  -5.46%  "
$ws.Range("D24").Value = "'76.74"
$ws.Range("E24").Value = "  -2.35%  "
$ws.Range("D25").Value = "'11.70"
$ws.Range("E25").Value = "  -2.85%  "
$ws.Range("D26").Value = "'0.993"
$ws.Range("E26").Value = "  -0.71%  "
$ws.Range("D27").Value = "'2.63"
$ws.Range("E27").Value = "  -2.47%  "
$ws.Range("D28").Value = "'7.52"
$ws.Range("E28").Value = "  -8.50%  "
$ws.Range("D29").Value = "'1.01"
$ws.Range("E29").Value = "  +1.09%  "
$ws.Range("D30").Value = "'25.19"
$ws.Range("E30").Value = "  -2.40%  "
$ws.Range("D31").Value = "'1.12"
$ws.Range("E31").Value = "  +1.10%  "
$ws.Range("D32").Value = "'1.80"
$ws.Range("E32").Value = "  -3.47%  "
$ws.Range("D33").Value = "'55.71"
$ws.Range("E33").Value = "  -3.70%  "
$ws.Range("D34").Value = "'2.21"
$ws.Range("E34").Value = "  -8.06%  "
$ws.Range("D35").Value = "'5.25"
$ws.Range("E35").Value = "  +0.70%  "
$ws.Range("D36").Value = "'5.72"
$ws.Range("E36").Value = "  -3.55%  "
$ws.Range("D37").Value = "'450.38"
$ws.Range("E37").Value = "  -4.72%  "
$ws.Range("D38").Value = "3.121.11"
$ws.Range("E38").Value = "  +0.25%  "
$ws.Range("D39").Value = "'0.0379"
$ws.Range("E39").Value = "  -2.92%  "
$ws.Range("D40").Value = "'0.0769"
$ws.Range("E40").Value = "  -2.61%  "
$ws.Range("D41").Value = "'0.115"
$ws.Range("E41").Value = "  +1.27%  "
$ws.Range("D42").Value = "'7.87"
$ws.Range("E42").Value = "  -2.10%  "
$ws.Range("D43").Value = "'2.39"
$ws.Range("E43").Value = "  -8.08%  "
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("D45").Value = "'0.241"
$ws.Range("E45").Value = "  -3.77%  "
$ws.Range("D46").Value = "'24.60"
$ws.Range("E46").Value = "  +0.98%  "
$ws.Range("D47").Value = "'120.14"
$ws.Range("E47").Value = "  +1.29%  "
$ws.Range("D48").Value = "'0.106"
$ws.Range("E48").Value = "  -0.75%  "
$ws.Range("D49").Value = "'1.89"
$ws.Range("E49").Value = "  -6.29%  "
$ws.Range("D50").Value = "0.0₃0495"
$ws.Range("E50").Value = "  -3.03%  "
$ws.Range("D51").Value = "'1.22"
$ws.Range("E51").Value = "  +4.52%  "
